$d = $word.ActiveDocument

# Near the end of the document, the references paragraph
# ("...GIL, A.C. ... 2010.") is followed by: an empty paragraph, the
# "Ver no Jupiter Salvar em pdf Salvar em docx" paragraph, and the
# "© 2020 ... Creative Commons Attribution" paragraph. All three of
# those need to be removed entirely (including their paragraph marks),
# so that the GIL paragraph is followed directly by the pre-existing
# empty paragraph that used to sit just before the page-break paragraph.

# Locate the "Ver no Jupiter ..." paragraph.
$findRng = $d.Content
$findRng.Find.Execute("Ver no Jupiter Salvar em pdf Salvar em docx", $true, $false, `
    $false, $false, $false, $true, 1, $false, "", 0)

# Find its index within Document.Paragraphs so we can reach its
# neighbours reliably (Paragraphs collections hanging off a Range don't
# support .Previous/.Next navigation in this host).
$allParas = $d.Paragraphs
$count = $allParas.Count
$targetIndex = -1
for ($i = 1; $i -le $count; $i++) {
    if ($allParas.Item($i).Range.Start -eq $findRng.Start) {
        $targetIndex = $i
    }
}

$startPara = $allParas.Item($targetIndex - 1)   # the empty paragraph right before it
$endPara = $allParas.Item($targetIndex + 1)      # the "© 2020 ..." paragraph right after it

# Delete the whole span (three paragraphs, including their paragraph
# marks) in one shot.
$deleteRange = $d.Range($startPara.Range.Start, $endPara.Range.End)
$deleteRange.Delete()
